# Add the two new "reference_period" and "remarks" fact fields as extra
# columns H and I on the import-file example sheet (row 1 headers only -
# the sample data row is not extended with values for these new columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the existing header/data range's cell style (re-applying the same
# "Normal" style) so it gets its own explicit style record, distinct from
# the default style the brand-new columns below will pick up. Restricted
# to the cells that already hold data so no new blank styled cell (G2) is
# materialised.
$ws.Range("A1:G1").Style = "Normal"
$ws.Range("A2:F2").Style = "Normal"

$ws.Range("H1").Value = "reference_period"
$ws.Range("I1").Value = "remarks"

# Give the new "reference_period" column a bit more breathing room, same
# as the other text columns.
$ws.Columns.Item(8).ColumnWidth = 14.73

# Leave the cursor where the user ended up after adding the columns.
$ws.Range("J11").Select() | Out-Null
